# Apply the edit described by the diff to the "metrics_12_9" worksheet.
#
# Effect observed in the diff:
#   1. Column A (model name labels) on rows 2-26 gets reordered according to
#      a fixed mapping (see below).
#   2. Columns B:Q on every data row (2-26) are overwritten with one single
#      repeated set of 16 numeric values (same set on every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New labels for column A, rows 2-26 -------------------------------
$newNames = @{
    2  = "model_12_9_0"
    3  = "model_12_9_22"
    4  = "model_12_9_21"
    5  = "model_12_9_20"
    6  = "model_12_9_19"
    7  = "model_12_9_18"
    8  = "model_12_9_17"
    9  = "model_12_9_16"
    10 = "model_12_9_15"
    11 = "model_12_9_14"
    12 = "model_12_9_13"
    13 = "model_12_9_23"
    14 = "model_12_9_12"
    15 = "model_12_9_10"
    16 = "model_12_9_9"
    17 = "model_12_9_8"
    18 = "model_12_9_7"
    19 = "model_12_9_6"
    20 = "model_12_9_5"
    21 = "model_12_9_4"
    22 = "model_12_9_3"
    23 = "model_12_9_2"
    24 = "model_12_9_1"
    25 = "model_12_9_11"
    26 = "model_12_9_24"
}

foreach ($row in $newNames.Keys) {
    $ws.Range("A$row").Value = $newNames[$row]
}

# --- 2. New constant values for columns B:Q, rows 2-26 -------------------
$newRowValues = @(
    0.6383931775788736,
    -23.10748969621655,
    0.5780077333516115,
    -0.5379569118765228,
    0.210832136944582,
    0.2146654303905574,
    14.31124727300474,
    0.1221825844847669,
    0.3496404292052329,
    0.2359115068449999,
    0.2718993657310753,
    0.4633200086231518,
    -0.08482046726337922,
    0.4830445005948765,
    35.07734920122535,
    54.57936239911656
)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $newRowValues[$i]
    }
}
